$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating "2022-Q1" (same columns,
#    headers and styling), positioned right after "总计" and before the
#    existing "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1    = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Helper scratch cell (far outside used range) used to push literal text
# values into cells without Excel's automatic "looks like a number -> make
# it a number" coercion, and without leaving any extra style/number-format
# behind on the destination cell.
$scratch = $wsQ3.Cells.Item(500, 26)

function Set-TextValue($cell, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

# Row 2 (090019 / 大成景恒混合A) — new 2022-Q3 figures.
Set-TextValue $wsQ3.Cells.Item(2, 4) "1.13"
Set-TextValue $wsQ3.Cells.Item(2, 5) "93.98"
Set-TextValue $wsQ3.Cells.Item(2, 6) "1.76"
Set-TextValue $wsQ3.Cells.Item(2, 7) "0.0199"
$wsQ3.Cells.Item(2, 8).Value = 8

# Row 3 (006038 / 大成景恒混合C) — new 2022-Q3 figures.
Set-TextValue $wsQ3.Cells.Item(3, 4) "0.45"
Set-TextValue $wsQ3.Cells.Item(3, 5) "93.98"
Set-TextValue $wsQ3.Cells.Item(3, 6) "1.76"
Set-TextValue $wsQ3.Cells.Item(3, 7) "0.0079"
$wsQ3.Cells.Item(3, 8).Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a row for 2022-Q3 right under
#    the header, pushing the existing "2022-Q1"/"2021-Q2" rows down one.
# ---------------------------------------------------------------------------
# A3 already carries the styled (s="2") index-column format; clone it onto
# the brand-new A4 before the old row-3 ("2021-Q2") data is moved into row 4.
$wsTotal.Cells.Item(3, 1).Copy($wsTotal.Cells.Item(4, 1)) | Out-Null

# Read old row 2/3 data first (.Value2 — plain .Value getter is unreliable
# in this host), then write the shifted rows.
$oldRow2B = $wsTotal.Cells.Item(2, 2).Value2
$oldRow2C = $wsTotal.Cells.Item(2, 3).Value2
$oldRow2D = $wsTotal.Cells.Item(2, 4).Value2
$oldRow3B = $wsTotal.Cells.Item(3, 2).Value2
$oldRow3C = $wsTotal.Cells.Item(3, 3).Value2
$oldRow3D = $wsTotal.Cells.Item(3, 4).Value2

# Shift old row 3 ("2021-Q2") down to row 4.
$wsTotal.Cells.Item(4, 2).Value = $oldRow3B
$wsTotal.Cells.Item(4, 3).Value = $oldRow3C
$wsTotal.Cells.Item(4, 4).Value = $oldRow3D

# Shift old row 2 ("2022-Q1") down to row 3.
$wsTotal.Cells.Item(3, 2).Value = $oldRow2B
$wsTotal.Cells.Item(3, 3).Value = $oldRow2C
$wsTotal.Cells.Item(3, 4).Value = $oldRow2D

# Row 2 becomes the new "2022-Q3" entry.
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 4).Value = 0.03

# Re-sequence the 0-based index column (A) for the three data rows.
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
